# Update "想去人数" (F) and "最低票价" (G) figures, and one cover-image URL (I7),
# on both the "展览" sheet and the "全部类型" sheet, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 342
    $ws.Range("F4").Value = 10560
    $ws.Range("F6").Value = 961

    $ws.Range("F7").Value = 84
    $ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202406/vF68pdMX1717579314139.png"

    $ws.Range("F9").Value = 7981
    $ws.Range("F10").Value = 24
    $ws.Range("F11").Value = 458

    $ws.Range("G12").Value = 60

    $ws.Range("F13").Value = 212
    $ws.Range("F15").Value = 3238
    $ws.Range("F17").Value = 321
    $ws.Range("F18").Value = 722
    $ws.Range("F20").Value = 1050
    $ws.Range("F21").Value = 282
    $ws.Range("F22").Value = 91
    $ws.Range("F23").Value = 1688
}
